$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (entity "98" / Utopia): ind_0114 (O), ind_0115 (P), ind_0116 (Q)
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = -1
$ws.Range("Q2").Value = -1

# Row 3 (entity "99" / Distopia): ind_0115 (P), ind_0116 (Q)
$ws.Range("P3").Value = -5
$ws.Range("Q3").Value = -5
